# Generate Report for Handoff
# b.md moved from "Handed back: in sync with en-US" to "Ready for handoff" with
# a fresh (newer) handoff package, and the handback version check now flags
# that the handback is stale relative to the new handoff.

$wb = $excel.ActiveWorkbook

$newStatus        = "Ready for handoff"
$newGenerateDate   = "2016-08-29 08:40:23"

$zhHandoffFile     = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate     = "2016-08-29 08:40:18"

$deHandoffFile     = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate     = "2016-08-29 08:40:23"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a3a33c9560ecb3280b64efe918e32082f81decd/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52f7a41be78d0283b67a1aae65952ba54b6b547b/e2e/b.md."

# --- Overview sheet: b.md row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newGenerateDate

# --- zh-cn sheet: b.md row (row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhHandoffFile
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Range("P1").ColumnWidth = 39.166666666666664

# --- de-de sheet: b.md row (row 3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deHandoffFile
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Range("P1").ColumnWidth = 39.166666666666664
